$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "Neaz "
$ws.Range("C5").Value = "Mahmud"

$ws.Range("C5").Select()
